$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 11, shifting existing rows 11-15 down to 12-16
$ws.Rows.Item(11).Insert()

# Give the new row the same formatting as the other feed-item rows (e.g. A10),
# then set its text.
$ws.Range("A10").Copy($ws.Range("A11"))
$excel.CutCopyMode = $false
$ws.Range("A11").Value = "Match Number in Match API"

# Update the selection to match the new state
$ws.Range("A1:A12").Select()
